$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-23T20:31:36+00:00"

# --- Include #0 sheet: prefix each Concept code with "NCIT:" ---
$inc = $wb.Worksheets.Item("Include #0")
$inc.Range("A2").Value = "NCIT:C156445"
$inc.Range("A3").Value = "NCIT:156440"
$inc.Range("A4").Value = "NCIT:156441"
$inc.Range("A5").Value = "NCIT:164032"
$inc.Range("A6").Value = "NCIT:C18009"
$inc.Range("A7").Value = "NCIT:C162623"
$inc.Range("A8").Value = "NCIT:C156443"
